$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.964.84'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").Value = '3.378.66'
$ws.Range("E3").Value = '  +7.45%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '260.32'
$ws.Range("E5").Value = '  +7.33%  '
$ws.Range("D6").Value = '629.60'
$ws.Range("E6").Value = '  +2.72%  '
$ws.Range("E7").Value = '  +23.54%  '
$ws.Range("D8").Value = '0.393'
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("B10").Value = 'LidoStakedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D10").Value = '3.377.06'
$ws.Range("E10").Value = '  +7.56%  '
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = '0.861'
$ws.Range("E11").Value = '  +10.10%  '
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '98.781.37'
$ws.Range("E13").Value = '  +2.13%  '
$ws.Range("D14").Value = '36.18'
$ws.Range("E14").Value = '  +5.63%  '
$ws.Range("E15").Value = '  +3.07%  '
$ws.Range("D16").Value = '4.005.26'
$ws.Range("E16").Value = '  +7.55%  '
$ws.Range("D17").Value = '5.50'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").Value = '3.380.85'
$ws.Range("E18").Value = '  +7.55%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '15.29'
$ws.Range("E20").Value = '  +4.32%  '
$ws.Range("D21").Value = '493.31'
$ws.Range("E21").Value = '  -4.88%  '
$ws.Range("E22").Value = '  +7.70%  '
$ws.Range("E23").Value = '  +8.67%  '
$ws.Range("E24").Value = '  +6.28%  '
$ws.Range("D25").Value = '5.64'
$ws.Range("E25").Value = '  +2.20%  '
$ws.Range("D26").Value = '88.56'
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").Value = '11.98'
$ws.Range("E27").Value = '  +2.57%  '
$ws.Range("D28").Value = '3.558.87'
$ws.Range("E28").Value = '  +7.76%  '
$ws.Range("E29").Value = '  +17.35%  '
$ws.Range("D31").Value = '0.192'
$ws.Range("E31").Value = '  +8.77%  '
$ws.Range("D32").Value = '0.130'
$ws.Range("E32").Value = '  +4.52%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '9.65'
$ws.Range("E33").Value = '  +6.43%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = '28.10'
$ws.Range("E35").Value = '  +5.11%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '7.33'
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.150'
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("E38").Value = '  +3.92%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '500.77'
$ws.Range("E39").Value = '  +4.09%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '0.461'
$ws.Range("E40").Value = '  +5.16%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = '24.91'
$ws.Range("E41").Value = '  +2.84%  '
$ws.Range("B42").Value = 'MantraDAO'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D42").Value = '3.83'
$ws.Range("E42").Value = '  +7.69%  '
$ws.Range("E43").Value = '  +2.56%  '
$ws.Range("E44").Value = '  +3.42%  '
$ws.Range("D45").Value = '0.788'
$ws.Range("E45").Value = '  +10.97%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '160.07'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").Value = '1.95'
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("D49").Value = '0.840'
$ws.Range("E49").Value = '  +12.97%  '
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("D51").Value = '46.07'
$ws.Range("E51").Value = '  +3.99%  '
